# Update the "Clasificación" standings table with the latest results,
# widen column A to fit the longer pair names, move the sheet selection,
# and make "Resultados" the active tab.

$wb = $excel.ActiveWorkbook

$wsClasificacion = $wb.Worksheets.Item("Clasificación")
$wsResultados    = $wb.Worksheets.Item("Resultados")

# --- New standings data (JJ, G, P, PTS(D), PTS) for rows 2-6 ---

# Pareja 26 - Jason/Jorge
$wsClasificacion.Range("B2").Value = 4
$wsClasificacion.Range("C2").Value = 2
$wsClasificacion.Range("D2").Value = 2
$wsClasificacion.Range("E2").Value = 367
$wsClasificacion.Range("F2").Value = 6

# Pareja 27 - Alex/Keneth
$wsClasificacion.Range("B3").Value = 4
$wsClasificacion.Range("C3").Value = 3
$wsClasificacion.Range("D3").Value = 1
$wsClasificacion.Range("E3").Value = 450
$wsClasificacion.Range("F3").Value = 9

# Pareja 28 - Teto/Pedro
$wsClasificacion.Range("B4").Value = 4
$wsClasificacion.Range("C4").Value = 1
$wsClasificacion.Range("D4").Value = 3
$wsClasificacion.Range("E4").Value = 290
$wsClasificacion.Range("F4").Value = 3

# Pareja 29 - Keko/Memo
$wsClasificacion.Range("B5").Value = 4
$wsClasificacion.Range("C5").Value = 2
$wsClasificacion.Range("D5").Value = 2
$wsClasificacion.Range("E5").Value = 310
$wsClasificacion.Range("F5").Value = 6

# Pareja 30 - Memin/Juan
$wsClasificacion.Range("B6").Value = 4
$wsClasificacion.Range("C6").Value = 3
$wsClasificacion.Range("D6").Value = 1
$wsClasificacion.Range("E6").Value = 360
$wsClasificacion.Range("F6").Value = 9

# Widen column A on "Clasificación" so the longer pair names fit
$wsClasificacion.Columns.Item(1).ColumnWidth = 14.666666666666668

# Move the selection on "Clasificación" to E6
$wsClasificacion.Range("E6").Select()

# Make "Resultados" the active (selected) tab
$wsResultados.Activate()
